$p = $ppt.ActivePresentation

# The edited text ("... -i input.i -o out" -> "... -i input.i -o out -text")
# lives in the "TextBox 18" shape on the "code interface" slide (slide 21).
$slide = $p.Slides.Item(21)

$target = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -like "*-o out*") {
            $target = $shape
            break
        }
    }
}

$tr = $target.TextFrame.TextRange
$tr.Replace(" -o out", " -o out -text", 0, $false, $false) | Out-Null
